$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the time-slot labels in column C for rows 7-11 to the new values.
$ws.Range("C7").Value2  = "22:30-22:35"
$ws.Range("C8").Value2  = "22:35-22:40"
$ws.Range("C9").Value2  = "22:40-22:45"
$ws.Range("C10").Value2 = "22:45-22:50"
$ws.Range("C11").Value2 = "22:50-22:55"

# Remove the now-obsolete last row (previously row 12: "22:25-22:30").
$ws.Rows.Item(12).Delete()

# Restore the previously-selected cell (now B15 instead of C15).
$ws.Range("B15").Select()
